$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9
$ws.Cells.Item($row - 1, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4122)
$ws.Cells.Item($row, 1).Value = 42612.892581018517
$ws.Cells.Item($row, 2).Value = 26
$ws.Cells.Item($row, 3).Value = 53
$ws.Cells.Item($row, 4).Value = 43
$ws.Cells.Item($row, 5).Value = 92
$ws.Cells.Item($row, 6).Value = 7
$ws.Cells.Item($row, 7).Value = 20669
$ws.Cells.Item($row, 8).Value = 16721
$ws.Cells.Item($row, 9).Value = 1027
$ws.Cells.Item($row, 10).Value = 174
$ws.Cells.Item($row, 11).Value = 140
$ws.Cells.Item($row, 12).Value = 13
$ws.Cells.Item($row, 13).Value = 1
$ws.Cells.Item($row, 14).Value = "Named"
